$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Court / probation / police.deja.demande"
$ws.Range("C1").Value = "General practitioner.deja.demande"
$ws.Range("D1").Value = "Other drug treatment centre.deja.demande"
$ws.Range("E1").Value = "Other health, medical or social service.deja.demande"
$ws.Range("F1").Value = "Educational services.deja.demande"
$ws.Range("G1").Value = "Self-referral, referral from family, friends, etc.; no other agency/institution involved.deja.demande"
$ws.Range("H1").Value = "Other.deja.demande"
$ws.Range("I1").Value = "Not known / missing.deja.demande"
$ws.Range("J1").Value = "Total.deja.demande"
